$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-09-16 Saturday" "2023-09-17 Sunday"

Replace-Text "52×21=1092" "76×61=4636"
Replace-Text "68×96=6528" "39×63=2457"
Replace-Text "31×71=2201" "80×68=5440"
Replace-Text "76×40=3040" "44×28=1232"
Replace-Text "83×62=5146" "65×54=3510"

Replace-Text "67×88=5896" "50×65=3250"
Replace-Text "61×67=4087" "56×77=4312"
Replace-Text "40×21=840" "30×43=1290"
Replace-Text "17×93=1581" "86×62=5332"
Replace-Text "15×94=1410" "97×83=8051"

Replace-Text "72×80=5760" "74×90=6660"
Replace-Text "57×53=3021" "53×51=2703"
Replace-Text "28×38=1064" "63×64=4032"
Replace-Text "75×75=5625" "93×67=6231"
Replace-Text "62×17=1054" "75×73=5475"

Replace-Text "63×82=5166" "15×76=1140"
Replace-Text "23×65=1495" "40×24=960"
Replace-Text "90×93=8370" "99×60=5940"
Replace-Text "84×21=1764" "53×19=1007"
Replace-Text "11×17=187" "27×95=2565"

Replace-Text "90×20=1800" "56×31=1736"
Replace-Text "93×29=2697" "16×34=544"
Replace-Text "75×95=7125" "45×17=765"
Replace-Text "94×35=3290" "85×33=2805"
Replace-Text "48×32=1536" "70×37=2590"
